$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.401.82"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.687.73"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'521.82"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "'146.22"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "2.705.58"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").Value = "'6.46"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D13").Value = "'0.129"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "3.150.28"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "60.422.96"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").Value = "'21.30"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "2.741.62"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'350.76"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").Value = "'4.55"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'10.55"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "'6.32"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'63.22"
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "0.0₃0814"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D30").Value = "'6.90"
$ws.Range("E30").Value = "  +7.53%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'19.07"
$ws.Range("D34").Value = "'148.68"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'4.35"
$ws.Range("E35").Value = "  +7.60%  "
$ws.Range("D36").Value = "'0.951"
$ws.Range("E36").Value = "  -6.03%  "
$ws.Range("E37").Value = "  +6.64%  "
$ws.Range("E38").Value = "  +11.07%  "
$ws.Range("D39").Value = "'0.873"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").Value = "'36.78"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'282.59"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'20.05"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0989"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.612"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "2.126.60"
$ws.Range("E47").Value = "  +7.31%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.85"
$ws.Range("E49").Value = "  +4.82%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0235"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.45"
$ws.Range("E51").Value = "  +2.63%  "
